# Sprint Planning workbook update: "Updated Sprint1 start date."
#  - Push the Sprint1 start date out by one week (2020-02-17 -> 2020-02-24)
#  - Reduce the API requirement from 3 APIs to 2 APIs in the baseline
#    requirements sheet (deliverable name + description)

$wb = $excel.ActiveWorkbook

# --- release sheet: move Sprint1 start date ---
$release = $wb.Worksheets.Item("release")
$release.Range("D2").Value = 43885
$release.Activate()
$release.Range("D3").Select()

# --- baselineReqs sheet: "Add 3 APIs" -> "Add 2 APIs" ---
$baseline = $wb.Worksheets.Item("baselineReqs")
$baseline.Range("A13").Value = "Add 2 APIs"
$baseline.Range("B13").Value = "Utilize 2 Web APIs to enhance functionality. Ex: payment processing, scheduling, mapping, social media"
$baseline.Activate()
$baseline.Range("B14").Select()
